$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.457.96'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.651.30'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'607.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').Value = "'156.07"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.15%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = '2.649.57'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('E10').Value = '  +7.67%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').Value = "'5.87"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').Value = "'29.92"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.39%  '
$ws.Range('D15').Value = "'0.0000196"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +15.07%  '
$ws.Range('D16').Value = '3.127.67'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').Value = '65.229.02'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('D18').Value = '2.656.26'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = "'12.73"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.48%  '
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').Value = "'358.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('D22').Value = "'7.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.69%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'69.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.14%  '
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = "'9.43"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('E27').Value = '  +14.92%  '
$ws.Range('D28').Value = "'1.63"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('E29').Value = '  +2.36%  '
$ws.Range('D30').Value = "'8.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.33%  '
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('E32').Value = '  +4.28%  '
$ws.Range('D33').Value = "'523.01"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.54%  '
$ws.Range('D34').Value = "'1.78"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('D35').Value = "'5.51"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Value = "'6.35"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('E37').Value = '  +2.22%  '
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('D39').Value = "'162.71"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.97"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = "'0.999"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D43').Value = "'41.89"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = "'165.47"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('E46').Value = '  +4.27%  '
$ws.Range('D47').Value = "'0.0610"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.27%  '
$ws.Range('D48').Value = "'22.95"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.650"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = "'0.0262"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.86%  '
$ws.Range('D51').Value = "'0.0979"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
